# Update Khoni disability_prevalence.xlsx:
#  - retitle the sheet header to the new "Unified database" description
#  - split the old single "Number of disability persons" row into two rows:
#      "family with disabilities Persons" and "disabilities Persons"
#    with freshly sourced data
#  - keep the existing "Source:" note (now on row 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the old data row (old row 4 -> 5, old row 5 -> 6)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Row 1 - title (merged A1:I1)
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Khoni Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").Font.ColorIndex = -4105
$ws.Range("A1:I1").Interior.Pattern = -4142
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A1:I1").Borders.Item(7).LineStyle = -4142
$ws.Range("A1:I1").Borders.Item(8).LineStyle = -4142
$ws.Range("A1:I1").Borders.Item(9).LineStyle = -4142
$ws.Range("A1:I1").Borders.Item(10).LineStyle = -4142
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Row 2 - "(End of year, persons)" caption (unchanged text/style, just
#    drop the explicit row height so Excel auto-sizes it again)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# 4. Row 3 - blank A3 cell switches font to Sylfaen 11 (years stay as-is)
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.ColorIndex = -4105
$ws.Range("A3").Font.Bold = $false

# ---------------------------------------------------------------------------
# 5. Row 4 (new) - "family with disabilities Persons " + data
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.ColorIndex = -4105
$ws.Range("A4").Font.Bold = $false
$ws.Range("A4").Interior.Pattern = -4142
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2

$row4vals = 929,940,938,957,900,866,809,800
$row4cols = "B","C","D","E","F","G","H","I"
for ($i = 0; $i -lt 8; $i++) {
    $addr = $row4cols[$i] + "4"
    $ws.Range($addr).Value = $row4vals[$i]
    $ws.Range($addr).NumberFormat = "#\ ##0"
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).Font.ColorIndex = -4105
    $ws.Range($addr).Interior.Pattern = -4142
}
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 6. Row 5 (old data row) - "disabilities Persons " + new data
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.ColorIndex = -4105
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").Interior.Pattern = -4142
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(7).LineStyle = 0
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").Borders.Item(10).LineStyle = 0

$row5vals = 1076,1073,1073,1085,1016,977,910,896
$row5cols = "B","C","D","E","F","G","H","I"
for ($i = 0; $i -lt 8; $i++) {
    $addr = $row5cols[$i] + "5"
    $ws.Range($addr).Value = $row5vals[$i]
    $ws.Range($addr).NumberFormat = "#\ ##0"
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).Font.ColorIndex = -4105
    $ws.Range($addr).Interior.Pattern = -4142
    $ws.Range($addr).Borders.Item(7).LineStyle = 0
    $ws.Range($addr).Borders.Item(8).LineStyle = 0
    $ws.Range($addr).Borders.Item(10).LineStyle = 0
}
# I5 additionally carries the bottom border of the table
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 7. Row 6 - "Source:" note, unchanged text, merged A6:H6
# ---------------------------------------------------------------------------
$ws.Range("A6:H6").Font.Name = "Arial"
$ws.Range("A6:H6").Font.Size = 9
$ws.Range("A6:H6").Font.ColorIndex = -4105
$ws.Range("A6:H6").Interior.Pattern = -4142
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true
$ws.Range("A6").Borders.Item(7).LineStyle = 0
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Range("A6").Borders.Item(9).LineStyle = 0
$ws.Range("A6").Borders.Item(10).LineStyle = 0
$ws.Range("B6:H6").Borders.Item(8).LineStyle = 1
$ws.Range("B6:H6").Borders.Item(8).Weight = 2
$ws.Range("A6:H6").Merge()
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column A width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20

# ---------------------------------------------------------------------------
# 9. Selection, matching the saved sheet view in the target file
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()

"done"
